$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.749.29"
$ws.Range("E2").Value = "  -1.46%  "

$ws.Range("D3").Value = "2.216.69"
$ws.Range("E3").Value = "  -1.36%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.91"
$ws.Range("E5").Value = "  +5.92%  "

$ws.Range("E6").Value = "  -0.52%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.33"
$ws.Range("E7").Value = "  +2.29%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.599"
$ws.Range("E9").Value = "  +7.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.61"
$ws.Range("E10").Value = "  +10.85%  "

$ws.Range("E11").Value = "  -3.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.22"
$ws.Range("E12").Value = "  -1.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.18"
$ws.Range("E13").Value = "  +6.13%  "

$ws.Range("E14").Value = "  -0.67%  "

$ws.Range("D15").Value = "2.545.66"
$ws.Range("E15").Value = "  -1.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.91"
$ws.Range("E16").Value = "  -1.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.869"
$ws.Range("E17").Value = "  -0.12%  "

$ws.Range("D18").Value = "2.210.34"
$ws.Range("E18").Value = "  -1.56%  "

$ws.Range("D19").Value = "41.674.23"
$ws.Range("E19").Value = "  -1.38%  "

$ws.Range("D20").Value = "0.0₃0960"
$ws.Range("E20").Value = "  -1.51%  "

$ws.Range("E21").Value = "  -1.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.70"
$ws.Range("E22").Value = "  -0.98%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.50"
$ws.Range("E23").Value = "  -0.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("E24").Value = "  +1.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.01"
$ws.Range("E25").Value = "  +9.25%  "

$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("E27").Value = "  +4.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.92"
$ws.Range("E28").Value = "  +8.80%  "

$ws.Range("E29").Value = "  -2.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.36"
$ws.Range("E30").Value = "  -0.34%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.71"
$ws.Range("E31").Value = "  +0.56%  "

$ws.Range("E32").Value = "  -1.29%  "

$ws.Range("E33").Value = "  +4.14%  "

$ws.Range("E34").Value = "  -2.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0736"
$ws.Range("E35").Value = "  +1.96%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.73"
$ws.Range("E36").Value = "  +0.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.46"
$ws.Range("E37").Value = "  +14.85%  "

$ws.Range("E38").Value = "  +6.27%  "

$ws.Range("E39").Value = "  +8.36%  "

$ws.Range("E40").Value = "  -0.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.93"
$ws.Range("E41").Value = "  -0.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.38"
$ws.Range("E42").Value = "  +20.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.45"
$ws.Range("E43").Value = "  -0.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.204"
$ws.Range("E44").Value = "  +6.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.81"
$ws.Range("E45").Value = "  -2.59%  "

$ws.Range("E46").Value = "  +2.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.67"
$ws.Range("E47").Value = "  -7.50%  "

$ws.Range("E48").Value = "  -1.58%  "

$ws.Range("E49").Value = "  -0.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.16"
$ws.Range("E50").Value = "  +4.40%  "

$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.18"
$ws.Range("E51").Value = "  -0.27%  "
